# Risikomanagement.xlsx: add the new "WPF Applikationen ... nicht flüssig" risk
# (row 11) to the "Risiken" table, and switch the active tab/selection from
# "Änderungsgeschichte" to "Risiken" (scrolled/selected around the new row).

$wb = $excel.ActiveWorkbook

$wsRisk = $wb.Worksheets.Item("Risiken")

# Activate the "Risiken" sheet -> becomes the workbook's active tab
# (workbookView activeTab="1") and picks up sheetView tabSelected="1",
# while "Änderungsgeschichte" loses its tabSelected flag.
$wsRisk.Activate()

# New risk entry (row 11): Nr. 8, risk title, description. The row grows
# tall enough (90pt) to show the wrapped description text.
$wsRisk.Rows.Item(11).RowHeight = 90
$wsRisk.Cells.Item(11, 1).Value = 8
$wsRisk.Cells.Item(11, 2).Value = "WPF Applikationen laufen bei hoher Auflösung nicht flüssig"
$wsRisk.Cells.Item(11, 3).Value = "Die Applikation besitzt zwar die Wunschauflösung, diese ist aber dadurch sind Animationen nicht mehr flüssig und die Applikation stürzt im schlimmsten Fall ab."

# Scroll so row 7 is at the top of the window, then select the new row's D
# cell, matching the saved view state (topLeftCell="A7", selection D11).
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$wsRisk.Range("D11").Select() | Out-Null

$wb.Save()
